# Apply cell updates per the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.369.09"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "3.371.86"
$ws.Range("E3").Value = "  -2.23%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'407.17"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'135.39"
$ws.Range("E6").Value = "  +10.21%  "

$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "'0.673"
$ws.Range("E9").Value = "  +3.47%  "

$ws.Range("E10").Value = "  -3.68%  "

$ws.Range("D11").Value = "'43.04"
$ws.Range("E11").Value = "  +4.25%  "

$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").Value = "3.903.80"
$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("D14").Value = "'8.37"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").Value = "'19.66"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").Value = "3.349.59"
$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("D17").Value = "61.359.38"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("D19").Value = "'11.00"
$ws.Range("E19").Value = "  +2.13%  "

$ws.Range("E20").Value = "  -4.29%  "

$ws.Range("D21").Value = "'3.21"
$ws.Range("E21").Value = "  -2.75%  "

$ws.Range("D22").Value = "'83.44"
$ws.Range("E22").Value = "  +2.48%  "

$ws.Range("D23").Value = "'314.42"
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("D24").Value = "'12.84"
$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").Value = "'4.77"
$ws.Range("E26").Value = "  +11.30%  "

$ws.Range("E27").Value = "  +7.41%  "

$ws.Range("D28").Value = "'29.44"
$ws.Range("E28").Value = "  -5.46%  "

$ws.Range("D29").Value = "'7.72"
$ws.Range("E29").Value = "  -2.11%  "

$ws.Range("D30").Value = "'0.117"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("D32").Value = "'11.33"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.49"
$ws.Range("E34").Value = "  -2.93%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'41.00"
$ws.Range("E35").Value = "  -2.33%  "

$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").Value = "'52.10"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  -2.40%  "

$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  -3.23%  "

$ws.Range("D41").Value = "'138.23"
$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.299"
$ws.Range("E42").Value = "  +6.24%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.97"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("D45").Value = "'4.05"
$ws.Range("E45").Value = "  +4.64%  "

$ws.Range("D46").Value = "'16.65"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("E47").Value = "  +1.59%  "

$ws.Range("D48").Value = "'21.49"
$ws.Range("E48").Value = "  -1.64%  "

$ws.Range("D49").Value = "2.130.27"
$ws.Range("E49").Value = "  -3.39%  "

$ws.Range("E50").Value = "  -4.99%  "

$ws.Range("D51").Value = "'1.91"
$ws.Range("E51").Value = "  +0.33%  "
